# Update weekly Fruta/Hortaliza data: dates (D), volumes (J) and prices (K,L,M,P)
# for rows 6 through 16 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 6;  D = 44259; J = 30; K = 4000; L = 4000; M = 4000; P = 4000 },
    @{ Row = 7;  D = 44313; J = 20; K = 4000; L = 4000; M = 4000; P = 4000 },
    @{ Row = 8;  D = 44176; J = 10; K = 4000; L = 4000; M = 4000; P = 4000 },
    @{ Row = 9;  D = 44291; J = 35; K = 4000; L = 4000; M = 4000; P = 4000 },
    @{ Row = 10; D = 44315; J = 40; K = 4000; L = 4000; M = 4000; P = 4000 },
    @{ Row = 11; D = 44316; J = 20; K = 4000; L = 4000; M = 4000; P = 4000 },
    @{ Row = 12; D = 44280; J = 55; K = 4000; L = 4000; M = 4000; P = 4000 },
    @{ Row = 13; D = 44508; J = 30; K = 4000; L = 4000; M = 4000; P = 4000 },
    @{ Row = 14; D = 44504; J = 55; K = 4000; L = 4000; M = 4000; P = 4000 },
    @{ Row = 15; D = 44301; J = 40; K = 3000; L = 3000; M = 3000; P = 3000 },
    @{ Row = 16; D = 44509; J = 20; K = 4000; L = 4000; M = 4000; P = 4000 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 4).Value  = $r.D   # D - Fecha
    $ws.Cells.Item($row, 10).Value = $r.J   # J - Volumen
    $ws.Cells.Item($row, 11).Value = $r.K   # K - Precio minimo
    $ws.Cells.Item($row, 12).Value = $r.L   # L - Precio maximo
    $ws.Cells.Item($row, 13).Value = $r.M   # M - Precio promedio ponderado
    $ws.Cells.Item($row, 16).Value = $r.P   # P - Precio $/Kg
}
